$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update animation duration reference values
$ws.Range("C4").Value = 0.8
$ws.Range("C5").Value = 1.3

# Set row 5 height (matches ht="15.75" customHeight="1")
$ws.Rows.Item(5).RowHeight = 15.75

# Update the selected cell / cursor position to D18
$ws.Range("D18").Select()
